# Insert two new rows right above the current row 930 (pushes existing
# rows 930-1048 down to 932-1050), then populate those two new rows with
# the new data points (one "Primera" and one "Segunda" reading dated
# serial 45124 = 2023-07-17).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(930).Insert()
$ws.Rows.Item(930).Insert()

# New row 930 - Primera
$ws.Cells.Item(930, 1).Value = 3
$ws.Cells.Item(930, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(930, 3).Value = "Coquimbo"
$ws.Cells.Item(930, 4).Value = 45124
$ws.Cells.Item(930, 5).Value = 5
$ws.Cells.Item(930, 6).Value = 100112006
$ws.Cells.Item(930, 7).Value = "Repollo"
$ws.Cells.Item(930, 8).Value = "Crespo record"
$ws.Cells.Item(930, 9).Value = "Primera"
$ws.Cells.Item(930, 10).Value = 9700
$ws.Cells.Item(930, 11).Value = 750
$ws.Cells.Item(930, 12).Value = 800
$ws.Cells.Item(930, 13).Value = 790
$ws.Cells.Item(930, 14).Value = "$/unidad"
$ws.Cells.Item(930, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(930, 16).Value = 790
$ws.Cells.Item(930, 17).Value = 1
$ws.Cells.Item(930, 18).Value = "Hortaliza"

# New row 931 - Segunda
$ws.Cells.Item(931, 1).Value = 3
$ws.Cells.Item(931, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(931, 3).Value = "Coquimbo"
$ws.Cells.Item(931, 4).Value = 45124
$ws.Cells.Item(931, 5).Value = 5
$ws.Cells.Item(931, 6).Value = 100112006
$ws.Cells.Item(931, 7).Value = "Repollo"
$ws.Cells.Item(931, 8).Value = "Crespo record"
$ws.Cells.Item(931, 9).Value = "Segunda"
$ws.Cells.Item(931, 10).Value = 1200
$ws.Cells.Item(931, 11).Value = 600
$ws.Cells.Item(931, 12).Value = 600
$ws.Cells.Item(931, 13).Value = 600
$ws.Cells.Item(931, 14).Value = "$/unidad"
$ws.Cells.Item(931, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(931, 16).Value = 600
$ws.Cells.Item(931, 17).Value = 1
$ws.Cells.Item(931, 18).Value = "Hortaliza"
